$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in row 34 with the new time-log entry
$ws.Range("B34").Value = 45355
$ws.Range("C34").Value = 1.5
$ws.Range("D34").Value = "Refactored search call method"

# Update the active selection to match the author's final cursor position
$ws.Range("D38").Select()
